$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 (week of 45943): speaker change + new note ---
$ws.Range("B17").Value = "Daniel"
$ws.Range("G17").Value = "Practice qual talk "

# --- Row 18 (week of 45950): clear speakers, update note to "No Lab Meeting" ---
$ws.Range("B18").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("G18").Value = "No Lab Meeting"

# --- Row 19 (week of 45957): speaker change + new note ---
$ws.Range("B19").Value = "Undergrads"
$ws.Range("G19").Value = "Joy, Giselle, Kevin.  Ben at AAO"

# --- Row 20 (week of 45964): add speaker + note ---
$ws.Range("B20").Value = "Everyone"
$ws.Range("G20").Value = "ARVO abstract talks"

# --- Row 21 (week of 45971): add speaker + note ---
$ws.Range("B21").Value = "Solomon"
$ws.Range("G21").Value = "Practice PhD dissertation"

# --- View: freeze top row, scroll down toward the bottom of the list, end selection at D21 ---
$win = $wb.Windows.Item(1)
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("D21").Select()
